$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 17.37

$ws.Range("B3").Value = 16.77
$ws.Range("C3").Value = 16.97

$ws.Range("C5").Value = 18.54
$ws.Range("C6").Value = 22.47
$ws.Range("C7").Value = 23.03
$ws.Range("C8").Value = 26.11
$ws.Range("C9").Value = 24.57
$ws.Range("C10").Value = 24.07
$ws.Range("C11").Value = 25.05
$ws.Range("C12").Value = 13.7
$ws.Range("C13").Value = 13.46
$ws.Range("C14").Value = 13.6

$ws.Range("C17").Value = 15.39
$ws.Range("C18").Value = 23.99
$ws.Range("C19").Value = 23.91
$ws.Range("C20").Value = 26.77
$ws.Range("C21").Value = 31.2
$ws.Range("C22").Value = 31.57

$ws.Range("C24").Value = 18.22
